$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "66.309.57"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "3.323.04"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.93"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.39"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +7.47%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "3.899.80"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("E13").Value = "  -4.12%  "
$ws.Range("D14").Value = "66.349.73"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.38"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "3.343.88"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "430.08"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.32"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.43"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.14"
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.70"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "3.460.25"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.203"
$ws.Range("E27").Value = "  +6.26%  "
$ws.Range("E28").Value = "  -4.09%  "
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.89"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").Value = "2.902.00"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.71"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.34"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.18"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0668"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.43"
$ws.Range("E48").Value = "  -5.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "316.74"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("E51").Value = "  +4.64%  "
